$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Treinamento")

# Fill in the missing sentiment-label scores (column B) for rows 136-200,
# continuing the manual labeling of the "Treinamento" sheet.
$pairs = @(
    @(136, 4),
    @(138, 2),
    @(140, 3),
    @(141, 3),
    @(142, 3),
    @(143, 3),
    @(144, 2),
    @(145, 4),
    @(146, 4),
    @(147, 2),
    @(148, 1),
    @(149, 2),
    @(150, 1),
    @(151, 2),
    @(152, 1),
    @(153, 4),
    @(154, 3),
    @(155, 1),
    @(156, 1),
    @(157, 3),
    @(158, 1),
    @(159, 4),
    @(160, 4),
    @(161, 4),
    @(163, 3),
    @(164, 1),
    @(165, 4),
    @(166, 0),
    @(167, 3),
    @(168, 1),
    @(169, 3),
    @(170, 1),
    @(171, 1),
    @(173, 3),
    @(174, 0),
    @(178, 0),
    @(179, 3),
    @(180, 1),
    @(181, 3),
    @(183, 2),
    @(184, 1),
    @(185, 2),
    @(186, 3),
    @(187, 4),
    @(188, 3),
    @(189, 2),
    @(190, 1),
    @(191, 2),
    @(192, 2),
    @(193, 1),
    @(194, 2),
    @(195, 3),
    @(197, 3),
    @(198, 4),
    @(199, 4),
    @(200, 4)
)

foreach ($pair in $pairs) {
    $r = $pair[0]
    $v = $pair[1]
    $ws.Cells.Item($r, 2).Value = $v
}

# Move the active selection to the next empty row, ready for further entry
$ws.Activate()
$ws.Range("B201").Select()
